$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "'25.54"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").Value = "'5.134"
$ws.Range("D4").ClearFormats()
$ws.Range("D6").Value = "'6.486"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'3.023"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'0.8178"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.8411"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.1335"
$ws.Range("D10").ClearFormats()
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.06961"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("D12").Value = "'0.02855"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'0.09385"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.001523"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'0.0005957"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '14OneONEWorstin24h'
$ws.Range("D16").Value = "'0.006169"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'3.532"
$ws.Range("D17").ClearFormats()
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = "'0.03210"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("D21").Value = "'0.1319"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'3.748"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").Value = "'0.04713"
$ws.Range("D23").ClearFormats()
$ws.Range("D25").Value = "'0.001246"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = "'0.004611"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").Value = "'0.00009697"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '26NitroExNTXBestin24h'
$ws.Range("D28").Value = "'0.0001389"
$ws.Range("D28").ClearFormats()
$ws.Range("D40").Value = "'0.03658"
$ws.Range("D40").ClearFormats()
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.006138"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1053"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = "'0.002529"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = "'0.007774"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.00005311"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").Value = "'0.1334"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").ClearFormats()
